$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "CanClone" column (K) — a single clone scene should no longer let the
# payer's login create/join an existing group; force these rows to 0 so a
# brand-new group gets created instead.
$ws.Range("K10").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("K15").Value = 0

# Reflect the new active selection left behind on the sheet view.
$ws.Range("K10").Select()
